$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character used in the PEPE price (U+2083).
# NOTE: build it as an explicit [string] and never use the bare
# [char] in a "+" concatenation -- this runtime treats "str" + [char]
# as numeric addition (e.g. "0.0" + [char]0x2083 becomes 8323.0).
$sub3 = [string][char]0x2083

# --- Column D (Price) / column E (Volume 1h) updates for rows whose
#     coin identity did not change. ---

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "68.050.02"
$ws.Cells.Item(2, 5).Value = "  -0.32%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.881.43"
$ws.Cells.Item(3, 5).Value = "  -1.19%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.17%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "481.95"
$ws.Cells.Item(5, 5).Value = "  -0.30%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "144.09"
$ws.Cells.Item(6, 5).Value = "  -2.50%  "

$ws.Cells.Item(7, 5).Value = "  -0.39%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.997"
$ws.Cells.Item(8, 5).Value = "  -0.06%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.738"
$ws.Cells.Item(9, 5).Value = "  +1.43%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.180"
$ws.Cells.Item(10, 5).Value = "  +7.80%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0000351"
$ws.Cells.Item(11, 5).Value = "  -0.49%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "42.72"
$ws.Cells.Item(12, 5).Value = "  +0.19%  "

$ws.Cells.Item(13, 5).Value = "  +0.95%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.497.65"

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.840.99"
$ws.Cells.Item(15, 5).Value = "  -3.09%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.24"
$ws.Cells.Item(16, 5).Value = "  -3.19%  "

$ws.Cells.Item(17, 5).Value = "  -0.66%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "20.03"
$ws.Cells.Item(18, 5).Value = "  +0.76%  "

$ws.Cells.Item(19, 5).Value = "  -0.57%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "68.065.79"
$ws.Cells.Item(20, 5).Value = "  -0.17%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "428.62"
$ws.Cells.Item(21, 5).Value = "  -1.47%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.57"
$ws.Cells.Item(22, 5).Value = "  +4.03%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "14.76"
$ws.Cells.Item(23, 5).Value = "  +1.96%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "89.92"
$ws.Cells.Item(24, 5).Value = "  +3.00%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "11.89"
$ws.Cells.Item(25, 5).Value = "  +9.25%  "

$ws.Cells.Item(26, 5).Value = "  +3.01%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.98"
$ws.Cells.Item(27, 5).Value = "  +2.77%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "37.36"
$ws.Cells.Item(28, 5).Value = "  -2.64%  "

$ws.Cells.Item(29, 5).Value = "  -3.83%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "711.19"
$ws.Cells.Item(30, 5).Value = "  -1.26%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "13.48"
$ws.Cells.Item(31, 5).Value = "  +1.43%  "

$ws.Cells.Item(32, 5).Value = "  +0.11%  "

$ws.Cells.Item(33, 5).Value = "  +2.55%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.07"
$ws.Cells.Item(34, 5).Value = "  +10.21%  "

$ws.Cells.Item(40, 5).Value = "  -3.98%  "

$ws.Cells.Item(42, 5).Value = "  +2.75%  "

$ws.Cells.Item(43, 5).Value = "  +3.37%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.97"
$ws.Cells.Item(44, 5).Value = "  -1.52%  "

$ws.Cells.Item(45, 5).Value = "  +0.78%  "

$ws.Cells.Item(46, 5).Value = "  +3.34%  "

$ws.Cells.Item(47, 5).Value = "  +0.29%  "

$ws.Cells.Item(48, 5).Value = "  -1.96%  "

$ws.Cells.Item(49, 5).Value = "  -3.24%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "144.16"
$ws.Cells.Item(50, 5).Value = "  -0.98%  "

$ws.Cells.Item(51, 5).Value = "  -1.80%  "

# Row 35 (PEPE) price contains a subscript-three character.
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0" + $sub3 + "0869"
$ws.Cells.Item(35, 5).Value = "  -2.27%  "

# --- Rows 36-41: ranking reshuffle. Coin identities (Coin/Link)
#     moved between rows; price and volume change accordingly. ---

$ws.Cells.Item(36, 2).Value = "OKB"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "60.88"
$ws.Cells.Item(36, 5).Value = "  +2.96%  "

$ws.Cells.Item(37, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "40.83"
$ws.Cells.Item(37, 5).Value = "  -2.76%  "

$ws.Cells.Item(38, 2).Value = "TheGraph"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.397"
$ws.Cells.Item(38, 5).Value = "  +14.02%  "

$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0500"
$ws.Cells.Item(39, 5).Value = "  +6.37%  "

$ws.Cells.Item(41, 2).Value = "Dai"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.997"
$ws.Cells.Item(41, 5).Value = "  -0.16%  "

